$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (index 0)
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = "9, 19, 31, 4, 11"
$ws.Range("D2").Value = "29 -> 31 -> 31 -> 19 -> 19 -> 11 -> 11 -> 4 -> 4 -> 9"
$ws.Range("E2").Value = 1188
$ws.Range("F2").Value = 0.1214172840118408

# Row 3 (index 1)
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = "10, 19, 8, 26, 13"
$ws.Range("D3").Value = "9 -> 13 -> 13 -> 19 -> 19 -> 10 -> 10 -> 8 -> 8 -> 26"
$ws.Range("E3").Value = 1172
$ws.Range("F3").Value = 0.1215567588806152

# Row 4 (index 2)
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = "6, 30, 21, 10, 4"
$ws.Range("D4").Value = "25 -> 6 -> 6 -> 4 -> 4 -> 10 -> 10 -> 21 -> 21 -> 30"
$ws.Range("E4").Value = 1163
$ws.Range("F4").Value = 0.1183443069458008

# Row 5 (index 3)
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "18, 24, 2, 8, 28"
$ws.Range("D5").Value = "5 -> 2 -> 2 -> 8 -> 8 -> 28 -> 28 -> 24 -> 24 -> 18"
$ws.Range("E5").Value = 1196
$ws.Range("F5").Value = 0.1190147399902344

# Row 6 (index 4)
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = "4, 17, 11, 18, 8"
$ws.Range("D6").Value = "8 -> 8 -> 8 -> 4 -> 4 -> 11 -> 11 -> 18 -> 18 -> 17"
$ws.Range("E6").Value = 869
$ws.Range("F6").Value = 0.09990334510803223

# Row 7 (index 5)
$ws.Range("B7").Value = 21
$ws.Range("C7").Value = "3, 4, 12, 23, 1"
$ws.Range("D7").Value = "21 -> 23 -> 23 -> 12 -> 12 -> 4 -> 4 -> 3 -> 3 -> 1"
$ws.Range("E7").Value = 1435
$ws.Range("F7").Value = 0.1154119968414307

# Row 8 (index 6)
$ws.Range("B8").Value = 25
$ws.Range("C8").Value = "11, 15, 6, 25, 9"
$ws.Range("D8").Value = "25 -> 25 -> 25 -> 6 -> 6 -> 9 -> 9 -> 11 -> 11 -> 15"
$ws.Range("E8").Value = 614
$ws.Range("F8").Value = 0.1045219898223877

# Row 9 (index 7)
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = "3, 30, 21, 7, 14"
$ws.Range("D9").Value = "6 -> 3 -> 3 -> 30 -> 30 -> 7 -> 7 -> 14 -> 14 -> 21"
$ws.Range("E9").Value = 1125
$ws.Range("F9").Value = 0.1311888694763184

# Row 10 (index 8)
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = "19, 10, 25, 26, 9"
$ws.Range("D10").Value = "19 -> 19 -> 19 -> 10 -> 10 -> 9 -> 9 -> 26 -> 26 -> 25"
$ws.Range("E10").Value = 848
$ws.Range("F10").Value = 0.1345524787902832

# Row 11 (index 9)
$ws.Range("B11").Value = 13
$ws.Range("C11").Value = "31, 9, 22, 28, 19"
$ws.Range("D11").Value = "13 -> 9 -> 9 -> 28 -> 28 -> 31 -> 31 -> 19 -> 19 -> 22"
$ws.Range("E11").Value = 836
$ws.Range("F11").Value = 0.1235864162445068
